$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UC")

# Update the estimated coding time for the "Export" use case (row 12, column C)
$ws.Range("C12").Value = 13

# Update the view to reflect the last worked-on area of the sheet
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("C13").Select()

$wb.Save()
